$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (tab 1): drop the f0d6173c... report row (row 3), update
# the handoff/handback report timestamps are NOT on this sheet, only the
# row removal + hyperlink cleanup applies here.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# Delete the stale report row (row 3) - shifts dimension down to A1:C2.
$wsOverview.Rows.Item(3).Delete()

# This engine's Range.Hyperlinks.Delete() clears the *whole sheet's*
# hyperlink collection, so clear once and rebuild only what should survive.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/adb180cab353f2e364aa07a0ff81d5a3d834cc22/e2e/ae18194b-63eb-4d66-9d9a-8155b62d5b35.md", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (tab 2): refresh the handoff/handback report timestamps for
# the surviving row, then drop the f0d6173c... report row (row 3).
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Range("E2").Value = "2016-03-17 06:07:24"
$wsZhCn.Range("H2").Value = "2016-03-17 06:08:03"

$wsZhCn.Rows.Item(3).Delete()

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/adb180cab353f2e364aa07a0ff81d5a3d834cc22/e2e/ae18194b-63eb-4d66-9d9a-8155b62d5b35.md", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/adb180cab353f2e364aa07a0ff81d5a3d834cc22/e2e/ae18194b-63eb-4d66-9d9a-8155b62d5b35.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09c8e38fd2639faf1d3282b9b66e2da178c68e4b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/0dba5000e20c035ec2865df716a9f891fe4104bb/e2e/ae18194b-63eb-4d66-9d9a-8155b62d5b35.md", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/368666c979b02208c79bac0fe3f4a914d23351af/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (tab 3): same treatment as zh-cn, with the de-de report
# timestamps/links.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Range("E2").Value = "2016-03-17 06:07:31"
$wsDeDe.Range("H2").Value = "2016-03-17 06:08:16"

$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/adb180cab353f2e364aa07a0ff81d5a3d834cc22/e2e/ae18194b-63eb-4d66-9d9a-8155b62d5b35.md", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/adb180cab353f2e364aa07a0ff81d5a3d834cc22/e2e/ae18194b-63eb-4d66-9d9a-8155b62d5b35.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a2c87e2994a47112da5a582883f401817462a73/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.de-de.xlf", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/b97c53c5aaca880dd1d6f4fa88e0a03c5f571490/e2e/ae18194b-63eb-4d66-9d9a-8155b62d5b35.md", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cdde3667c9d75f05de0fee05cf054c33a870aa92/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.de-de.xlf", [Type]::Missing, [Type]::Missing, "ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.de-de.xlf") | Out-Null
